$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 1..6 + 8..11 + 13..17 + 19..22 + 24..27
foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = 1
}

$ws.Range("D28").Select()
